# Nederlands Profiel voor Catalogi - content update
# (title rename StelselCatalogi -> Catalogi, new explanatory texts for the
#  "Concept Schema's", "Concept verzamelingen" and "Overige eigenschappen"
#  sections, removal of the stray "bibliographic resource" note, and the
#  scrolled/selected view state left behind by the author.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Begrip")
$ws.Activate()

# Title: "Profiel voor StelselCatalogi" -> "Profiel voor Catalogi"
$ws.Range("A1").Value = "Profiel voor Catalogi"

# Section intro for "Concept Schema's (begrippenkaders)" (merged D21:L21)
$ws.Range("D21").Value = "het Concept Schema (in het Nederlands begrippenkader) is de Universe of Discourse, en verzorgt de functionele invulling om context vast te leggen. "

# Section intro for "Concept verzamelingen (Collections)" (merged D23:L23)
$ws.Range("D23").Value = "gelabelde of geordende groepen van SKOS concepten, bruikbaar om concepten die een overeenkomst hebben te grouperen en te labellen.  "

# Section intro for "Overige eigenschappen" (merged D25:L25)
$ws.Range("D25").Value = "Aanvullende eigenschappen welke niet direct afkomstig zijn uit SKOS."

# Drop the leftover "bibliographic resource" note under "Bron"
$ws.Range("D26").ClearContents()

# Restore the scroll/selection state the author left the sheet in:
# frozen pane still splits after row 3, but scrolled down to row 19,
# with A25 selected.
$pane2 = $excel.ActiveWindow.Panes.Item(2)
$pane2.Activate()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A25").Select()

$wb.Save()
